$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Year_selection")

# Final target state for the Year/Selected-Years table (A2:B14).
# Row 6 used to be 2025/1 -> becomes 2021/0
# Rows 7-10 shift the existing year values down one row (2030->2025, 2035->2030, 2040->2035, 2045->2040)
# Row 11 becomes 2045/1 (new), row 12 becomes 2050/1 (the old row 11 content), rows 13-14 are brand new (2055/1, 2060/1)
$years = @(2021, 2025, 2030, 2035, 2040, 2045, 2050, 2055, 2060)
$flags = @(0, 1, 1, 1, 1, 1, 1, 1, 1)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = 6 + $i
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $flags[$i]
}

# Update the active selection to A7 (matches the saved sheetView state in the target file)
[void]$ws.Range("A7").Select()
